# update data March 24
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("China_Sales_Vehicle_YTD")

# Correct the previously-provisional row 243 (2023-12-31) with final figures
$ws.Range("C243").Value = 9495236
$ws.Range("D243").Value = 9047861
$ws.Range("E243").Value = 6685176
$ws.Range("F243").Value = 2804217
$ws.Range("G243").Value = 5843

# Add new row 244 (2024-01-31) - copy formatting from row 243 first
$ws.Range("A243:H243").Copy()
$ws.Range("A244:H244").PasteSpecial(-4122)

$ws.Range("A244").Value = 45322
$ws.Range("B244").Value = 243.9
$ws.Range("C244").Value = 729317
$ws.Range("D244").Value = 698901
$ws.Range("E244").Value = 444647
$ws.Range("F244").Value = 284292
$ws.Range("G244").Value = 378
$ws.Range("H244").Value = 29.9

# Add new row 245 (2024-02-29, provisional) - copy formatting from row 243
$ws.Range("A243:H243").Copy()
$ws.Range("A245:H245").PasteSpecial(-4122)

$ws.Range("A245").Value = 45351
$ws.Range("B245").Value = 402.6
$ws.Range("C245").Value = 1207000
$ws.Range("D245").Value = 0
$ws.Range("E245").Value = 740000
$ws.Range("F245").Value = 467000
$ws.Range("G245").Value = 1000
$ws.Range("H245").Value = 29.98
